# Updates cryptos list values (price/volume columns, plus two coin-rank swaps)
# Commit: "Updated cryptos list on Wed Mar 15 18:49:22 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'24.471.29"
$ws.Range("E2").Value = "  -4.27%  "

# Row 3
$ws.Range("D3").Value = "'1.644.53"
$ws.Range("E3").Value = "  -6.21%  "

# Row 4
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.52%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'306.04"
$ws.Range("E5").Value = "  -3.45%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'0.9979"
$ws.Range("E6").Value = "  -0.38%  "

# Row 7
$ws.Range("D7").Value = "'0.3622"
$ws.Range("E7").Value = "  -5.24%  "

# Row 8
$ws.Range("D8").Value = "'47.31"
$ws.Range("E8").Value = "  -4.57%  "

# Row 9
$ws.Range("D9").Value = "'0.3257"
$ws.Range("E9").Value = "  -9.18%  "

# Row 10
$ws.Range("D10").Value = "'1.118"
$ws.Range("E10").Value = "  -8.71%  "

# Row 11
$ws.Range("D11").Value = "'0.06907"
$ws.Range("E11").Value = "  -9.97%  "

# Row 12
$ws.Range("D12").Value = "'0.9969"
$ws.Range("E12").Value = "  -0.41%  "

# Row 13
$ws.Range("D13").Value = "'5.923"
$ws.Range("E13").Value = "  -8.67%  "

# Row 14
$ws.Range("D14").Value = "'19.10"
$ws.Range("E14").Value = "  -11.29%  "

# Row 15
$ws.Range("D15").Value = "'1.646.09"
$ws.Range("E15").Value = "  -6.24%  "

# Row 16
$ws.Range("D16").Value = "'6.522"
$ws.Range("E16").Value = "  -8.50%  "

# Row 17
$ws.Range("D17").Value = "'0.00001044"
$ws.Range("E17").Value = "  -9.65%  "

# Row 18
$ws.Range("D18").Value = "'0.06487"
$ws.Range("E18").Value = "  -4.25%  "

# Row 19
$ws.Range("D19").Value = "'0.9979"
$ws.Range("E19").Value = "  -0.41%  "

# Row 20
$ws.Range("D20").Value = "'76.78"
$ws.Range("E20").Value = "  -11.01%  "

# Row 21
$ws.Range("D21").Value = "'5.883"
$ws.Range("E21").Value = "  -9.93%  "

# Row 22
$ws.Range("D22").Value = "'15.70"
$ws.Range("E22").Value = "  -11.00%  "

# Row 23
$ws.Range("D23").Value = "'12.14"
$ws.Range("E23").Value = "  -7.86%  "

# Row 24
$ws.Range("D24").Value = "'24.404.81"
$ws.Range("E24").Value = "  -4.40%  "

# Row 25
$ws.Range("D25").Value = "'2.396"
$ws.Range("E25").Value = "  -2.20%  "

# Row 26
$ws.Range("D26").Value = "'2.326"
$ws.Range("E26").Value = "  -19.82%  "

# Row 27
$ws.Range("D27").Value = "'145.09"
$ws.Range("E27").Value = "  -6.12%  "

# Row 28
$ws.Range("E28").Value = "  -11.67%  "

# Row 29
$ws.Range("D29").Value = "'1.825.32"
$ws.Range("E29").Value = "  -6.35%  "

# Row 30
$ws.Range("D30").Value = "'124.51"
$ws.Range("E30").Value = "  -7.08%  "

# Row 31
$ws.Range("D31").Value = "'1.142"
$ws.Range("E31").Value = "  -5.62%  "

# Row 32
$ws.Range("D32").Value = "'4.064"
$ws.Range("E32").Value = "  -3.62%  "

# Row 33
$ws.Range("D33").Value = "'5.599"
$ws.Range("E33").Value = "  -22.01%  "

# Row 34
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'1.685"
$ws.Range("E34").Value = "  -6.89%  "

# Row 35
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.08312"
$ws.Range("E35").Value = "  -5.05%  "

# Row 36
$ws.Range("D36").Value = "'12.34"
$ws.Range("E36").Value = "  -14.14%  "

# Row 37
$ws.Range("D37").Value = "'5.135"
$ws.Range("E37").Value = "  -11.19%  "

# Row 38
$ws.Range("D38").Value = "'0.06038"
$ws.Range("E38").Value = "  -10.10%  "

# Row 39
$ws.Range("D39").Value = "'0.02210"
$ws.Range("E39").Value = "  -11.51%  "

# Row 40
$ws.Range("D40").Value = "'8.199"
$ws.Range("E40").Value = "  -12.82%  "

# Row 41
$ws.Range("D41").Value = "'1.199"
$ws.Range("E41").Value = "  -7.22%  "

# Row 42
$ws.Range("D42").Value = "'0.2035"
$ws.Range("E42").Value = "  -10.12%  "

# Row 43
$ws.Range("D43").Value = "'0.9966"
$ws.Range("E43").Value = "  -0.45%  "

# Row 44
$ws.Range("D44").Value = "'0.5838"
$ws.Range("E44").Value = "  -11.44%  "

# Row 45
$ws.Range("D45").Value = "'3.720"
$ws.Range("E45").Value = "  -4.75%  "

# Row 46
$ws.Range("D46").Value = "'12.52"
$ws.Range("E46").Value = "  -13.00%  "

# Row 47
$ws.Range("D47").Value = "'0.5583"
$ws.Range("E47").Value = "  -11.44%  "

# Row 48
$ws.Range("D48").Value = "'121.41"
$ws.Range("E48").Value = "  -8.01%  "

# Row 49
$ws.Range("D49").Value = "'1.929"
$ws.Range("E49").Value = "  -11.87%  "

# Row 50
$ws.Range("D50").Value = "'0.06895"
$ws.Range("E50").Value = "  -7.30%  "

# Row 51
$ws.Range("D51").Value = "'73.53"
$ws.Range("E51").Value = "  -9.37%  "
